$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J, copy H1's style+value then overwrite text
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 4

$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 4
